# Correct the treatment column on the ZEV biosample sheet: the cells that
# were labeled "EtOH" are actually the mock-induction (no real estradiol)
# control and should read "mockEstradiol" instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$used.Replace("EtOH", "mockEstradiol", 2)
